$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9407587051391602
$ws.Range("B1").Value = 1.532688498497009
$ws.Range("C1").Value = 5.892971515655518
$ws.Range("D1").Value = 1.741750717163086
$ws.Range("E1").Value = 1.071824789047241
